$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# Header for new column D
$ws.Range("D1").Value = "Cost"

# Fill D2:D23 with 1
$ws.Range("D2:D23").Value = 1

# Select D2:D23 with D2 as the active cell, and activate the Tasks sheet/tab
$ws.Activate()
$ws.Range("D2:D23").Select()
